# Generate Report for Handback
# Update the localization-status report with the results of a failed
# handback transform for the 8db4a7bf-... file in both target languages.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Every place that used to read "Ready for handoff" for the 8db4a7bf-...
# file (row 3) now reads "Handback transform failed" - the Overview
# summary columns for each language plus the per-language Status column.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C3").Value     = $newStatus
$dede.Range("C3").Value     = $newStatus

# zh-cn sheet: record the handback/handoff file-name mismatch error in the
# "Error Detail" column (P) for row 3, and widen that column so the message
# is readable (the COM layer's character-width -> XML-width round trip adds
# ~0.83, so back the input off slightly to land on exactly 40).
$zhcn.Range("P3").Value = "Handback file name: 4pm3l5ad.b5l is different with handoff file name: 8db4a7bf-7179-4a72-aa11-8596bd6a207f.ea94d06a0d339b23488b7d95aca0cb9245218d08.zh-cn."
$zhcn.Columns.Item(16).ColumnWidth = 39.15

# de-de sheet: same error for the German handback.
$dede.Range("P3").Value = "Handback file name: 4pm3l5ad.b5l is different with handoff file name: 8db4a7bf-7179-4a72-aa11-8596bd6a207f.ea94d06a0d339b23488b7d95aca0cb9245218d08.de-de."
$dede.Columns.Item(16).ColumnWidth = 39.15
